$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-07-09 Tuesday" "2024-07-10 Wednesday"

Replace-Text "282×3=" "588×9="
Replace-Text "296×7=" "337×9="
Replace-Text "569×7=" "816×5="
Replace-Text "867×9=" "926×3="
Replace-Text "378×2=" "324×3="
Replace-Text "478×2=" "898×7="
Replace-Text "923×2=" "926×5="
Replace-Text "825×3=" "823×4="
Replace-Text "710×2=" "468×5="
Replace-Text "145×8=" "658×9="
Replace-Text "525×8=" "531×9="
Replace-Text "954×9=" "811×6="
Replace-Text "710×9=" "205×9="
Replace-Text "905×6=" "831×8="
Replace-Text "434×6=" "605×3="
Replace-Text "224×2=" "986×9="
Replace-Text "936×4=" "910×4="
Replace-Text "784×6=" "778×9="
Replace-Text "668×7=" "417×2="
Replace-Text "399×2=" "325×7="
Replace-Text "900×4=" "598×3="
Replace-Text "640×2=" "987×9="
Replace-Text "356×9=" "818×5="
Replace-Text "250×5=" "540×9="
Replace-Text "563×2=" "319×3="
